$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 232, pushing existing rows 232-252 down to 233-253.
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new weekly record.
$ws.Cells.Item(232, 1).Value = 9
$ws.Cells.Item(232, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(232, 3).Value = "Metropolitana"
$ws.Cells.Item(232, 4).Value = 44783
$ws.Cells.Item(232, 5).Value = 13
$ws.Cells.Item(232, 6).Value = 100112026
$ws.Cells.Item(232, 7).Value = "Haba"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 65
$ws.Cells.Item(232, 11).Value = 13000
$ws.Cells.Item(232, 12).Value = 14000
$ws.Cells.Item(232, 13).Value = 13462
$ws.Cells.Item(232, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(232, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(232, 16).Value = 538
$ws.Cells.Item(232, 17).Value = 25
$ws.Cells.Item(232, 18).Value = "Hortaliza"
